# Add a new "category" column between "body" (B) and "when" (C), then
# populate the header and the per-row category values. The existing
# "when" column shifts from C to D automatically as part of the insert.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at C; this shifts the old C ("when") to D and
# carries the header style (bold/border/centered) along with it.
$ws.Columns.Item(3).Insert()

# New header for the inserted column.
$ws.Range("C1").Value = "category"

# Shorten the body text for row 2 and set its category.
$ws.Range("B2").Value = @"
## 촬영 전 확인
연령, 매복치 여부 등
## 이미지 포인트
- 치낭낭종: 균일, 얇은 피질
- OKC: scalloping, 재발률↑
"@
$ws.Range("C2").Value = "영상감별"

# Shorten the body text for row 3 and set its category.
$ws.Range("B3").Value = @"
## 임상과 영상 결합
- 농양: 급성, 경계 불명확
- 낭종: 만성, 경계 명확
"@
$ws.Range("C3").Value = "근단병소"
